# Inserts a new weekly record at row 53 of the "Perejil" sheet, shifting all
# subsequent rows (old 53-112) down by one (to 54-113), and fills the newly
# inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 53 - this shifts rows 53..112
# down to 54..113 and extends the sheet dimension automatically.
$ws.Rows(53).Insert()

# Populate the newly inserted row 53 with the new record's values.
$ws.Range("A53").Value = 7
$ws.Range("B53").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C53").Value = "Ñuble"
$ws.Range("D53").Value = 45159
$ws.Range("E53").Value = 16
$ws.Range("F53").Value = 100112044
$ws.Range("G53").Value = "Perejil"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 200
$ws.Range("K53").Value = 1500
$ws.Range("L53").Value = 1500
$ws.Range("M53").Value = 1500
$ws.Range("N53").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O53").Value = "Región de Ñuble"
$ws.Range("P53").Value = 1500
$ws.Range("Q53").Value = 1
$ws.Range("R53").Value = "Hortaliza"
